$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.794230699539185
$ws.Range("B1").Value = 1.932799935340881
$ws.Range("C1").Value = 1.889144420623779
$ws.Range("D1").Value = 2.258776426315308
$ws.Range("E1").Value = 3.180292129516602
